{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items');\nawait context.sync();\n\n// Paragraph 0: replace text, then add a line break + second line\nparagraphs.items[0].insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 26.06.25\", Word.InsertLocation.replace);\nparagraphs.items[0].insertBreak(Word.BreakType.line, Word.InsertLocation.end);\nparagraphs.items[0].insertText(\"Open Problems in Mechanistic Interpretability\", Word.InsertLocation.end);\nawait context.sync();\n\n// Paragraphs 1-11: replace text in place\nconst middleTexts = [\n  \"\u05d0\u05d9\u05e0\u05d8\u05e8\u05e4\u05e8\u05d8\u05d1\u05d9\u05dc\u05d9\u05d5\u05ea \u05de\u05db\u05e0\u05d9\u05e1\u05d8\u05d9\u05ea \u05d4\u05d9\u05d0 \u05d0\u05d5\u05dc\u05d9 \u05ea\u05d7\u05d5\u05dd \u05d4\u05e9\u05d0\u05e4\u05ea\u05e0\u05d9 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05db\u05d9\u05d5\u05dd \u05dc\u05d4\u05d1\u05e0\u05ea \u05d0\u05d9\u05da \u05d1\u05d9\u05e0\u05d4 \u05de\u05dc\u05d0\u05db\u05d5\u05ea\u05d9\u05ea \u05d1\u05d0\u05de\u05ea \u05e2\u05d5\u05d1\u05d3\u05ea. \u05dc\u05d0 \u05de\u05d3\u05d5\u05d1\u05e8 \u05db\u05d0\u05df \u05d1\u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05d1\u05e0\u05e4\u05e0\u05d5\u05e4\u05d9 \u05d9\u05d3\u05d9\u05d9\u05dd \u05d0\u05d5 \u05d1\u05d4\u05d3\u05d2\u05e9\u05d5\u05ea \u05e6\u05d1\u05e2\u05d5\u05e0\u05d9\u05d5\u05ea \u05e9\u05dc \u05d7\u05dc\u05e7\u05d9 \u05d8\u05e7\u05e1\u05d8 \u05d0\u05dc\u05d0 \u05d1\u05d4\u05e0\u05d3\u05e1\u05d4 \u05dc\u05d0\u05d7\u05d5\u05e8(reverse engineering) \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea\u05d5\u05ea \u05e2\u05e6\u05de\u05df. \u05d4\u05d1\u05e0\u05d4 \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05e9\u05dc \u05d0\u05d9\u05da \u05e8\u05e9\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05e4\u05d5\u05ea\u05e8\u05ea \u05d1\u05e2\u05d9\u05d4: \u05de\u05d4\u05dd \u05d4\u05d7\u05dc\u05e7\u05d9\u05dd \u05d4\u05e4\u05e0\u05d9\u05de\u05d9\u05d9\u05dd \u05e9\u05e4\u05d5\u05e2\u05dc\u05d9\u05dd, \u05d1\u05d0\u05d9\u05d6\u05d4 \u05e1\u05d3\u05e8, \u05d1\u05d0\u05d9\u05d6\u05d5 \u05dc\u05d5\u05d2\u05d9\u05e7\u05d4, \u05d5\u05d0\u05d9\u05da \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d4\u05dd \u05de\u05d9\u05d9\u05e6\u05e8\u05d9\u05dd \u05d4\u05db\u05dc\u05dc\u05d4. \",\n  \"\u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05de\u05e8\u05db\u05d6\u05d9\u05ea \u05e9\u05de\u05d5\u05e6\u05d2\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc \u05e9\u05dc\u05d5\u05e9\u05d4 \u05e9\u05dc\u05d1\u05d9\u05dd: \u05e4\u05d9\u05e8\u05d5\u05e7 \u05d4\u05e8\u05e9\u05ea \u05dc\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e7\u05d8\u05e0\u05d9\u05dd (\u05d1\u05d9\u05df \u05d0\u05dd \u05d0\u05dc\u05d5 \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd, \u05ea\u05ea\u05d9-\u05de\u05e8\u05d7\u05d1\u05d9\u05dd \u05d0\u05d5 \u05de\u05e2\u05d2\u05dc\u05d9\u05dd), \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d4\u05ea\u05e4\u05e7\u05d9\u05d3 \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9 \u05e9\u05dc \u05db\u05dc \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd, \u05d5\u05d0\u05d9\u05de\u05d5\u05ea \u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05d3\u05d9\u05e7\u05d4 \u05d4\u05d0\u05dd \u05d4\u05d4\u05e1\u05d1\u05e8 \u05e9\u05dc\u05e0\u05d5 \u05d1\u05d0\u05de\u05ea \u05d7\u05d5\u05d6\u05d4 \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea, \u05d5\u05d0\u05dd \u05db\u05df \u05e2\u05d3 \u05db\u05de\u05d4. \u05db\u05dc \u05d0\u05d7\u05d3 \u05de\u05d4\u05e9\u05dc\u05d1\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05de\u05ea\u05d2\u05dc\u05d4 \u05db\u05e7\u05e9\u05d4 \u05d4\u05e8\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de\u05de\u05d4 \u05e9\u05e0\u05d3\u05de\u05d4.\",\n  \"\u05d4\u05d1\u05e2\u05d9\u05d4 \u05d4\u05d1\u05e1\u05d9\u05e1\u05d9\u05ea \u05d4\u05d9\u05d0 \u05e9\u05e4\u05d9\u05e8\u05d5\u05e7 \u05dc\u05e4\u05d9 \u05de\u05d1\u05e0\u05d4 \u05d4\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea \u05db\u05dc\u05d5\u05de\u05e8 \u05e9\u05db\u05d1\u05d5\u05ea, \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd, \u05e8\u05d0\u05e9\u05d9 attention \u05e4\u05e9\u05d5\u05d8 \u05dc\u05d0 \u05e2\u05d5\u05d1\u05d3. \u05d4\u05d7\u05dc\u05e7\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05dc\u05d0 \u05de\u05ea\u05d0\u05d9\u05de\u05d9\u05dd \u05dc\u05de\u05d4 \u05e9\u05d4\u05e8\u05e9\u05ea \u05d1\u05d0\u05de\u05ea \u05de\u05d7\u05e9\u05d1\u05ea. \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05d4\u05dd \u05e4\u05d5\u05dc\u05d9\u05e1\u05de\u05e0\u05d8\u05d9\u05d9\u05dd(\u05e8\u05d1 \u05de\u05e9\u05de\u05e2\u05d9\u05dd), \u05ea\u05e4\u05e7\u05d9\u05d3\u05d9\u05dd \u05de\u05ea\u05e4\u05e8\u05e9\u05d9\u05dd \u05e2\u05dc \u05e4\u05e0\u05d9 \u05e9\u05db\u05d1\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea, \u05d5\u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05dc\u05d0 \u05e9\u05d5\u05db\u05e0\u05d5\u05ea \u05d1\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d1\u05d5\u05d3\u05d3 \u05d0\u05dc\u05d0 \u05de\u05e7\u05d5\u05d3\u05d3\u05d5\u05ea \u05db\u05e1\u05d5\u05e4\u05e8\u05e4\u05d5\u05d6\u05d9\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05e8\u05d1\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05d4\u05e7\u05dc\u05d0\u05e1\u05d9\u05d5\u05ea \u05db\u05de\u05d5 PCA \u05d5-SVD \u05e0\u05db\u05e9\u05dc\u05d5\u05ea, \u05dc\u05d0 \u05d1\u05d2\u05dc\u05dc \u05d9\u05d9\u05e9\u05d5\u05dd \u05dc\u05e7\u05d5\u05d9 \u05d0\u05dc\u05d0 \u05d1\u05d2\u05dc\u05dc \u05d4\u05e0\u05d7\u05d5\u05ea \u05ea\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05e9\u05d2\u05d5\u05d9\u05d5\u05ea.\",\n  \"\u05d4\u05db\u05dc\u05d9 \u05d4\u05de\u05e8\u05db\u05d6\u05d9 \u05db\u05d9\u05d5\u05dd \u05d4\u05d5\u05d0 Sparse Dictionary Learning \u05d5\u05d1\u05e2\u05d9\u05e7\u05e8 Sparse Autoencoders. \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05d5\u05d0 \u05dc\u05d0\u05de\u05df \u05e8\u05e9\u05ea \u05e7\u05d8\u05e0\u05d4 \u05e9\\\"\u05ea\u05e4\u05e8\u05e9\\\" \u05d0\u05ea \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea \u05d4\u05d2\u05d3\u05d5\u05dc\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d1\u05e1\u05d9\u05e1 \u05d3\u05dc\u05d9\u05dc \u05e9\u05dc \\\"\u05ea\u05db\u05d5\u05e0\u05d5\u05ea\\\". \u05d0\u05dc\u05d5 \u05d4\u05dc\u05d9\u05d9\u05d8\u05e0\u05d8\u05d9\u05dd. \u05d0\u05da \u05d1\u05e4\u05d5\u05e2\u05dc, \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d0\u05de\u05e0\u05dd \u05de\u05d5\u05e6\u05d0\u05ea \u05db\u05d9\u05d5\u05d5\u05e0\u05d9\u05dd \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d9\u05dd, \u05d0\u05da \u05dc\u05d0 \u05de\u05e1\u05d1\u05d9\u05e8\u05d4 \u05d0\u05d9\u05da \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05e2\u05e6\u05de\u05d5 \u05de\u05ea\u05d1\u05e6\u05e2. \u05d4\u05dc\u05d9\u05d9\u05d8\u05e0\u05d8\u05d9\u05dd \u05d4\u05dd \u05ea\u05de\u05d5\u05e0\u05d4 \u05e1\u05d8\u05d8\u05d9\u05ea \u05e9\u05dc \\\"\u05de\u05d4 \u05d4\u05d5\u05e4\u05e2\u05dc\\\" \u05d5\u05dc\u05d0 \u05ea\u05d9\u05d0\u05d5\u05e8 \u05e9\u05dc \u05d4\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05e9\u05de\u05d9\u05d5\u05e9\u05dd.\",\n  \"\u05d9\u05e9 \u05d2\u05dd \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05d4\u05d5\u05ea\u05d9\u05d5\u05ea: \u05d4\u05e4\u05e2\u05e8 \u05d1\u05d9\u05df \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d4\u05d0\u05de\u05d9\u05ea\u05d9\u05d5\u05ea \u05dc\u05e9\u05d9\u05d7\u05d6\u05d5\u05e8\u05df \u05d2\u05d3\u05d5\u05dc. \u05d4\u05de\u05d9\u05d3\u05e2 \u05d4\u05d2\u05d0\u05d5\u05de\u05d8\u05e8\u05d9 \u05d1\u05d9\u05df \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d4\u05d5\u05dc\u05da \u05dc\u05d0\u05d9\u05d1\u05d5\u05d3. \u05d4\u05d4\u05e0\u05d7\u05d4 \u05e9\u05d4\u05db\u05d5\u05dc \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 \u05e8\u05d7\u05d5\u05e7\u05d4 \u05de\u05dc\u05d4\u05d9\u05d5\u05ea \u05e0\u05db\u05d5\u05e0\u05d4. \u05d5\u05d4\u05d2\u05e8\u05d5\u05e2 \u05de\u05db\u05dc \u05d4\u05d9\u05d0 \u05d4\u05e2\u05d5\u05d1\u05d3\u05d4 \u05e9\u05d0\u05d9\u05df \u05d1\u05db\u05dc\u05dc \u05ea\u05d9\u05d0\u05d5\u05e8\u05d9\u05d4 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e9\u05de\u05e1\u05d1\u05d9\u05e8\u05d4 \u05de\u05d4\u05d9 \\\"\u05ea\u05db\u05d5\u05e0\u05d4\\\", \u05d0\u05d9\u05da \u05d4\u05d9\u05d0 \u05e0\u05d5\u05e6\u05e8\u05ea, \u05d5\u05de\u05d4 \u05d4\u05d5\u05e4\u05da \u05d0\u05d5\u05ea\u05d4 \u05dc\u05d9\u05d7\u05d9\u05d3\u05d4 \u05d1\u05e1\u05d9\u05e1\u05d9\u05ea \u05e9\u05dc \u05d4\u05d1\u05e0\u05d4.\",\n  \"\u05de\u05db\u05d0\u05df \u05e2\u05d5\u05dc\u05d4 \u05db\u05d9\u05d5\u05d5\u05df \u05d7\u05d3\u05e9\u05e0\u05d9: \u05d0\u05d5\u05dc\u05d9 \u05d4\u05d3\u05e8\u05da \u05d4\u05e0\u05db\u05d5\u05e0\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d0 \u05dc\u05e4\u05e8\u05e9 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05d7\u05e8\u05d9 \u05e9\u05d0\u05d5\u05de\u05e0\u05d5, \u05d0\u05dc\u05d0 \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05e8\u05d0\u05e9. \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e2\u05dd \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d3\u05d9\u05e1\u05e7\u05e8\u05d8\u05d9\u05d5\u05ea, \u05d0\u05db\u05d9\u05e4\u05ea \u05de\u05d5\u05d3\u05d5\u05dc\u05e8\u05d9\u05d5\u05ea, \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05e4\u05e2\u05dc\u05d4 \u05d3\u05dc\u05d9\u05dc\u05d5\u05ea \u05db\u05de\u05d5 Top-k \u05d0\u05d5 SoLU, \u05d0\u05d5 \u05de\u05d1\u05e0\u05d9\u05dd \u05db\u05de\u05d5 Mixture-of-Experts \u05e9\u05de\u05d7\u05dc\u05e7\u05d9\u05dd \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05dc\u05ea\u05ea-\u05de\u05d5\u05d3\u05d5\u05dc\u05d9\u05dd \u05d1\u05e8\u05d5\u05e8\u05d9\u05dd. \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d9\u05d9\u05e6\u05e8 \u05e8\u05e9\u05ea\u05d5\u05ea \u05e9\u05e0\u05d1\u05e0\u05d5\u05ea \\\"\u05d7\u05ea\u05d5\u05db\u05d5\u05ea \u05de\u05e8\u05d0\u05e9\\\" \u05e2\u05dd \u05e4\u05e8\u05e9\u05e0\u05d5\u05ea \u05dc\u05d0 \u05db\u05e0\u05d9\u05ea\u05d5\u05d7 \u05de\u05d0\u05d5\u05d7\u05e8 \u05d0\u05dc\u05d0 \u05db\u05d4\u05e0\u05d7\u05ea \u05d9\u05e1\u05d5\u05d3 \u05e9\u05dc \u05d4\u05d0\u05d9\u05de\u05d5\u05df.\",\n  \"\u05d2\u05dd \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e8\u05db\u05d9\u05d1 \u05d1\u05d5\u05d3\u05d3 \u05d4\u05d5\u05d0 \u05de\u05e9\u05d9\u05de\u05d4 \u05e7\u05e9\u05d4. \u05dc\u05de\u05e9\u05dc \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05de\u05e4\u05e2\u05d9\u05dc\u05d5\u05ea \u05d0\u05d5\u05ea\u05d5 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d1\u05dc\u05d1\u05dc\u05d5\u05ea. \u05e9\u05d9\u05d8\u05d5\u05ea \u05d9\u05d9\u05d7\u05d5\u05e1 \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05d1\u05e2\u05d9\u05d9\u05ea\u05d9\u05d5\u05ea \u05ea\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d5\u05e4\u05e8\u05e7\u05d8\u05d9\u05ea. \u05e1\u05d9\u05e0\u05ea\u05d6\u05ea \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d1\u05e0\u05d9\u05d9\u05ea \u05e7\u05dc\u05d8 \u05e9\u05de\u05e4\u05e2\u05d9\u05dc \u05e8\u05db\u05d9\u05d1 \u05e2\u05dc\u05d5\u05dc\u05d4 \u05dc\u05d9\u05d9\u05e6\u05e8 \u05d3\u05d9\u05de\u05d5\u05d9\u05d9\u05dd \u05dc\u05d0 \u05d0\u05d9\u05e0\u05e4\u05d5\u05e8\u05de\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05d4\u05de\u05d1\u05d8\u05d9\u05d7\u05d5\u05ea \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05df \u05d0\u05d9\u05e0\u05d8\u05e8\u05d5\u05d5\u05e0\u05e6\u05d9\u05d5\u05ea \u05e1\u05d9\u05d1\u05ea\u05d9\u05d5\u05ea: \u05e9\u05d9\u05e0\u05d5\u05d9 \u05e9\u05dc \u05e2\u05e8\u05da \u05e4\u05e0\u05d9\u05de\u05d9, \u05d5\u05d1\u05d7\u05d9\u05e0\u05d4 \u05e9\u05dc \u05d4\u05d4\u05e9\u05e4\u05e2\u05d4 \u05e2\u05dc \u05d4\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05d7\u05d9\u05e6\u05d5\u05e0\u05d9\u05ea. \u05db\u05d0\u05df \u05e0\u05db\u05e0\u05e1\u05d9\u05dd \u05dc\u05ea\u05de\u05d5\u05e0\u05d4 \u05d2\u05dd steering \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d7\u05d3\u05e8\u05d4 \u05e9\u05dc \u05db\u05d9\u05d5\u05d5\u05df \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9 \u05dc\u05de\u05e8\u05d7\u05d1 \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d5\u05d2\u05dd \u05e9\u05d9\u05de\u05d5\u05e9 logit lens \u05db\u05d3\u05d9 \u05dc\u05e4\u05e2\u05e0\u05d7 \u05d4\u05e9\u05e4\u05e2\u05d4 \u05d9\u05e9\u05d9\u05e8\u05d4 \u05e2\u05dc \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05e2\u05dc \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d4\u05e8\u05e9\u05ea.\",\n  \"\u05d4\u05d1\u05e2\u05d9\u05d4 \u05d4\u05d2\u05d3\u05d5\u05dc\u05d4 \u05d4\u05d9\u05d0 \u05e9\u05d4\u05e8\u05d1\u05d4 \u05de\u05d4\u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05e0\u05e9\u05de\u05e2\u05d9\u05dd \u05de\u05e9\u05db\u05e0\u05e2\u05d9\u05dd \u05d0\u05da \u05dc\u05d0 \u05e2\u05d5\u05de\u05d3\u05d9\u05dd \u05d1\u05de\u05d1\u05d7\u05df. \u05d4\u05dd \u05dc\u05d0 \u05d7\u05d5\u05d6\u05d9\u05dd \u05e7\u05d5\u05e0\u05d8\u05e8\u05e4\u05e7\u05d8\u05d5\u05d0\u05dc\u05d9\u05dd(\u05dc\u05d0 \u05de\u05e6\u05dc\u05d9\u05d7\u05d9\u05dd \u05dc\u05e0\u05d1\u05d0 \u05de\u05d4 \u05d4\u05d9\u05d4 \u05e7\u05d5\u05e8\u05d4 \u05d0\u05d9\u05dc\u05d5 \u05de\u05e9\u05d4\u05d5 \u05d4\u05d9\u05d4 \u05e9\u05d5\u05e0\u05d4 \u05d1\u05ea\u05d5\u05da \u05d4\u05de\u05d5\u05d3\u05dc), \u05dc\u05d0 \u05e2\u05d5\u05d6\u05e8\u05d9\u05dd \u05dc\u05d0\u05d1\u05d7\u05df \u05db\u05e9\u05dc\u05d9 \u05de\u05d5\u05d3\u05dc, \u05dc\u05d0 \u05de\u05d0\u05e4\u05e9\u05e8\u05d9\u05dd \u05ea\u05d9\u05e7\u05d5\u05df \u05d0\u05d5 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d1\u05e4\u05d5\u05e2\u05dc. \u05dc\u05db\u05df \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e1\u05d8 \u05e9\u05dc\u05dd \u05e9\u05dc \u05d3\u05e8\u05db\u05d9 \u05d0\u05d9\u05de\u05d5\u05ea: \u05d4\u05d0\u05dd \u05d4\u05d4\u05e1\u05d1\u05e8 \u05d7\u05d5\u05d6\u05d4 \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d0\u05d7\u05e8\u05d9 ablation? \u05d4\u05d0\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc \u05e7\u05d8\u05df \u05e9\u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05dd \u05d4\u05d4\u05e1\u05d1\u05e8 \u05e0\u05db\u05d5\u05df? \u05d4\u05d0\u05dd \u05d4\u05dc\u05d9\u05d9\u05d8\u05e0\u05d8\u05d9\u05dd ( \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e4\u05e0\u05d9\u05de\u05d9\u05d9\u05dd \u05db\u05de\u05d5 \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d0\u05d5 \u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05d9\u05dd) \u05de\u05e1\u05d9\u05d9\u05e2\u05d9\u05dd \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d1\u05d8\u05d9\u05d7\u05d5\u05ea \u05db\u05de\u05d5 \u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d5\u05db\u05df \u05de\u05d6\u05d9\u05e7? \u05d4\u05d0\u05dd \u05e0\u05d5\u05db\u05dc \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05e9\u05e0\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05de\u05d5\u05d3\u05dc?\",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05d2\u05dd \u05dc\u05d9\u05e6\u05d5\u05e8 \\\"\u05d0\u05d5\u05e8\u05d2\u05e0\u05d9\u05d6\u05de\u05d9\u05dd \u05de\u05d5\u05d3\u05dc\u05d9\u05d9\u05dd\\\" \u05e9\u05d4\u05dd \u05de\u05d4\u05d5\u05d5\u05d9\u05dd \u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d8\u05e0\u05d5\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05d5\u05ea, \u05e2\u05dd \u05de\u05d1\u05e0\u05d4 \u05e4\u05ea\u05d5\u05d7, \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d0\u05de\u05df \u05e9\u05d5\u05d1 \u05d5\u05e9\u05d5\u05d1 \u05d5\u05dc\u05d1\u05d3\u05d5\u05e7 \u05e2\u05dc\u05d9\u05d4\u05df \u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05e8\u05e9\u05e0\u05d5\u05ea. \u05db\u05de\u05d5 \u05e9\u05d4\u05d1\u05d9\u05d5\u05dc\u05d5\u05d2\u05d9\u05d4 \u05d4\u05ea\u05e7\u05d3\u05de\u05d4 \u05d3\u05e8\u05da \u05e2\u05d1\u05d5\u05d3\u05d4 \u05e2\u05dc \u05ea\u05e1\u05d9\u05e1\u05e0\u05d9\u05ea, \u05db\u05da \u05ea\u05d7\u05d5\u05dd \u05d6\u05d4 \u05d6\u05e7\u05d5\u05e7 \u05dc\u05e8\u05e4\u05e8\u05e0\u05e1 \u05e7\u05d1\u05d5\u05e2. \u05d6\u05d4\u05d5 \u05db\u05dc\u05d9 \u05ea\u05e9\u05ea\u05d9\u05ea\u05d9 \u05d7\u05e1\u05e8 \u05db\u05d9\u05d5\u05dd.\",\n  \"\u05d4\u05d7\u05dc\u05e7 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d1\u05d4\u05d9\u05e8 \u05e9\u05de\u05db\u05e0\u05d9\u05d6\u05dd \u05d0\u05d9\u05e0\u05d5 \u05e2\u05e0\u05d9\u05d9\u05df \u05d8\u05db\u05e0\u05d9 \u05d1\u05dc\u05d1\u05d3. \u05d4\u05d5\u05d0 \u05e0\u05d5\u05d2\u05e2 \u05dc\u05de\u05d3\u05d9\u05e0\u05d9\u05d5\u05ea, \u05dc\u05e0\u05d9\u05d8\u05d5\u05e8, \u05dc\u05d1\u05d8\u05d9\u05d7\u05d5\u05ea, \u05d5\u05dc\u05e9\u05d0\u05dc\u05d5\u05ea \u05e4\u05d9\u05dc\u05d5\u05e1\u05d5\u05e4\u05d9\u05d5\u05ea: \u05de\u05d4 \u05e0\u05d7\u05e9\u05d1 \u05d4\u05e1\u05d1\u05e8 \u05d8\u05d5\u05d1? \u05d0\u05d9\u05da \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d7\u05d1\u05e8 \u05d1\u05d9\u05df \u05d4\u05de\u05d1\u05e0\u05d9\u05dd \u05d4\u05de\u05d9\u05e7\u05e8\u05d5\u05e1\u05e7\u05d5\u05e4\u05d9\u05d9\u05dd \u05dc\u05ea\u05e4\u05e7\u05d5\u05d3 \u05d2\u05dc\u05d5\u05d1\u05dc\u05d9? \u05d0\u05d9\u05dc\u05d5 \u05e2\u05e7\u05e8\u05d5\u05e0\u05d5\u05ea \u05db\u05dc\u05dc\u05d9\u05d9\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05d7\u05dc\u05e5 \u05de\u05e8\u05e9\u05ea\u05d5\u05ea \u05e9\u05dc\u05de\u05d3\u05d5 \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05d8\u05d5\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05de\u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd?\",\n  \"\u05d1\u05e1\u05d9\u05db\u05d5\u05dd, \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05de\u05ea\u05d1\u05d9\u05d9\u05e9 \u05dc\u05d5\u05de\u05e8 \u05d0\u05ea \u05d4\u05d0\u05de\u05ea: \u05d0\u05d9\u05df \u05e2\u05d3\u05d9\u05d9\u05df \u05ea\u05d9\u05d0\u05d5\u05e8\u05d9\u05d4 \u05de\u05e1\u05e4\u05e7\u05ea \u05dc\u05e4\u05d9\u05e8\u05d5\u05e7 \u05e8\u05e9\u05ea\u05d5\u05ea. \u05d4\u05d4\u05e0\u05d7\u05d5\u05ea \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d5\u05ea \u05e9\u05d1\u05e8\u05d9\u05e8\u05d9\u05d5\u05ea. \u05d4\u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05dc\u05d0 \u05d7\u05d9\u05d5\u05ea \u05dc\u05d1\u05d3 \u05d0\u05dc\u05d0 \u05d1\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05e2\u05dc. \u05d4\u05e4\u05e8\u05e9\u05e0\u05d5\u05ea \u05d7\u05d9\u05d9\u05d1\u05ea \u05dc\u05e7\u05e9\u05d5\u05e8 \u05de\u05d1\u05e0\u05d4 \u05dc\u05ea\u05e4\u05e7\u05d5\u05d3. \u05d5\u05d4\u05d3\u05e8\u05da \u05e7\u05d3\u05d9\u05de\u05d4, \u05d0\u05d5\u05dc\u05d9, \u05e2\u05d5\u05d1\u05e8\u05ea \u05dc\u05d0 \u05d3\u05e8\u05da \u05e4\u05e2\u05e0\u05d5\u05d7 \u05d0\u05dc\u05d0 \u05d3\u05e8\u05da \u05f4\u05d3\u05d6\u05d9\u05d9\u05df\u05f4 \u05d7\u05d3\u05e9 \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea\u05d5\u05ea\u2026\",\n];\nfor (let i = 0; i < middleTexts.length; i++) {\n  paragraphs.items[i + 1].insertText(middleTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Paragraphs 12-18 (old indices): delete entirely (7 paragraphs)\nfor (let i = 0; i < 7; i++) {\n  paragraphs.items[12].delete();\n}\nawait context.sync();\n\n// Last paragraph: replace URL text\nparagraphs.load('items');\nawait context.sync();\nparagraphs.items[paragraphs.items.length - 1].insertText(\"https://arxiv.org/abs/2501.16496\", Word.InsertLocation.replace);\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 1: replace text, insert a line break, then second line\n$p1 = $d.Paragraphs.Item(1)\n$p1.Range.Text = '\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 26.06.25'\n$endRange = $p1.Range\n$endRange.Collapse(0)  # wdCollapseEnd\n$endRange.InsertBreak(6)  # wdLineBreak\n$endRange.InsertAfter('Open Problems in Mechanistic Interpretability')\n\n# Paragraphs 2-12: replace text in place\n$middleTexts = @(\n  '\u05d0\u05d9\u05e0\u05d8\u05e8\u05e4\u05e8\u05d8\u05d1\u05d9\u05dc\u05d9\u05d5\u05ea \u05de\u05db\u05e0\u05d9\u05e1\u05d8\u05d9\u05ea \u05d4\u05d9\u05d0 \u05d0\u05d5\u05dc\u05d9 \u05ea\u05d7\u05d5\u05dd \u05d4\u05e9\u05d0\u05e4\u05ea\u05e0\u05d9 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05db\u05d9\u05d5\u05dd \u05dc\u05d4\u05d1\u05e0\u05ea \u05d0\u05d9\u05da \u05d1\u05d9\u05e0\u05d4 \u05de\u05dc\u05d0\u05db\u05d5\u05ea\u05d9\u05ea \u05d1\u05d0\u05de\u05ea \u05e2\u05d5\u05d1\u05d3\u05ea. \u05dc\u05d0 \u05de\u05d3\u05d5\u05d1\u05e8 \u05db\u05d0\u05df \u05d1\u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05d1\u05e0\u05e4\u05e0\u05d5\u05e4\u05d9 \u05d9\u05d3\u05d9\u05d9\u05dd \u05d0\u05d5 \u05d1\u05d4\u05d3\u05d2\u05e9\u05d5\u05ea \u05e6\u05d1\u05e2\u05d5\u05e0\u05d9\u05d5\u05ea \u05e9\u05dc \u05d7\u05dc\u05e7\u05d9 \u05d8\u05e7\u05e1\u05d8 \u05d0\u05dc\u05d0 \u05d1\u05d4\u05e0\u05d3\u05e1\u05d4 \u05dc\u05d0\u05d7\u05d5\u05e8(reverse engineering) \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea\u05d5\u05ea \u05e2\u05e6\u05de\u05df. \u05d4\u05d1\u05e0\u05d4 \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05e9\u05dc \u05d0\u05d9\u05da \u05e8\u05e9\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05e4\u05d5\u05ea\u05e8\u05ea \u05d1\u05e2\u05d9\u05d4: \u05de\u05d4\u05dd \u05d4\u05d7\u05dc\u05e7\u05d9\u05dd \u05d4\u05e4\u05e0\u05d9\u05de\u05d9\u05d9\u05dd \u05e9\u05e4\u05d5\u05e2\u05dc\u05d9\u05dd, \u05d1\u05d0\u05d9\u05d6\u05d4 \u05e1\u05d3\u05e8, \u05d1\u05d0\u05d9\u05d6\u05d5 \u05dc\u05d5\u05d2\u05d9\u05e7\u05d4, \u05d5\u05d0\u05d9\u05da \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d4\u05dd \u05de\u05d9\u05d9\u05e6\u05e8\u05d9\u05dd \u05d4\u05db\u05dc\u05dc\u05d4. ',\n  '\u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05de\u05e8\u05db\u05d6\u05d9\u05ea \u05e9\u05de\u05d5\u05e6\u05d2\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc \u05e9\u05dc\u05d5\u05e9\u05d4 \u05e9\u05dc\u05d1\u05d9\u05dd: \u05e4\u05d9\u05e8\u05d5\u05e7 \u05d4\u05e8\u05e9\u05ea \u05dc\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e7\u05d8\u05e0\u05d9\u05dd (\u05d1\u05d9\u05df \u05d0\u05dd \u05d0\u05dc\u05d5 \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd, \u05ea\u05ea\u05d9-\u05de\u05e8\u05d7\u05d1\u05d9\u05dd \u05d0\u05d5 \u05de\u05e2\u05d2\u05dc\u05d9\u05dd), \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d4\u05ea\u05e4\u05e7\u05d9\u05d3 \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9 \u05e9\u05dc \u05db\u05dc \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd, \u05d5\u05d0\u05d9\u05de\u05d5\u05ea \u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05d3\u05d9\u05e7\u05d4 \u05d4\u05d0\u05dd \u05d4\u05d4\u05e1\u05d1\u05e8 \u05e9\u05dc\u05e0\u05d5 \u05d1\u05d0\u05de\u05ea \u05d7\u05d5\u05d6\u05d4 \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea, \u05d5\u05d0\u05dd \u05db\u05df \u05e2\u05d3 \u05db\u05de\u05d4. \u05db\u05dc \u05d0\u05d7\u05d3 \u05de\u05d4\u05e9\u05dc\u05d1\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05de\u05ea\u05d2\u05dc\u05d4 \u05db\u05e7\u05e9\u05d4 \u05d4\u05e8\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de\u05de\u05d4 \u05e9\u05e0\u05d3\u05de\u05d4.',\n  '\u05d4\u05d1\u05e2\u05d9\u05d4 \u05d4\u05d1\u05e1\u05d9\u05e1\u05d9\u05ea \u05d4\u05d9\u05d0 \u05e9\u05e4\u05d9\u05e8\u05d5\u05e7 \u05dc\u05e4\u05d9 \u05de\u05d1\u05e0\u05d4 \u05d4\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea \u05db\u05dc\u05d5\u05de\u05e8 \u05e9\u05db\u05d1\u05d5\u05ea, \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd, \u05e8\u05d0\u05e9\u05d9 attention \u05e4\u05e9\u05d5\u05d8 \u05dc\u05d0 \u05e2\u05d5\u05d1\u05d3. \u05d4\u05d7\u05dc\u05e7\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05dc\u05d0 \u05de\u05ea\u05d0\u05d9\u05de\u05d9\u05dd \u05dc\u05de\u05d4 \u05e9\u05d4\u05e8\u05e9\u05ea \u05d1\u05d0\u05de\u05ea \u05de\u05d7\u05e9\u05d1\u05ea. \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05d4\u05dd \u05e4\u05d5\u05dc\u05d9\u05e1\u05de\u05e0\u05d8\u05d9\u05d9\u05dd(\u05e8\u05d1 \u05de\u05e9\u05de\u05e2\u05d9\u05dd), \u05ea\u05e4\u05e7\u05d9\u05d3\u05d9\u05dd \u05de\u05ea\u05e4\u05e8\u05e9\u05d9\u05dd \u05e2\u05dc \u05e4\u05e0\u05d9 \u05e9\u05db\u05d1\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea, \u05d5\u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05dc\u05d0 \u05e9\u05d5\u05db\u05e0\u05d5\u05ea \u05d1\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d1\u05d5\u05d3\u05d3 \u05d0\u05dc\u05d0 \u05de\u05e7\u05d5\u05d3\u05d3\u05d5\u05ea \u05db\u05e1\u05d5\u05e4\u05e8\u05e4\u05d5\u05d6\u05d9\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05e8\u05d1\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05d4\u05e7\u05dc\u05d0\u05e1\u05d9\u05d5\u05ea \u05db\u05de\u05d5 PCA \u05d5-SVD \u05e0\u05db\u05e9\u05dc\u05d5\u05ea, \u05dc\u05d0 \u05d1\u05d2\u05dc\u05dc \u05d9\u05d9\u05e9\u05d5\u05dd \u05dc\u05e7\u05d5\u05d9 \u05d0\u05dc\u05d0 \u05d1\u05d2\u05dc\u05dc \u05d4\u05e0\u05d7\u05d5\u05ea \u05ea\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea \u05e9\u05d2\u05d5\u05d9\u05d5\u05ea.',\n  '\u05d4\u05db\u05dc\u05d9 \u05d4\u05de\u05e8\u05db\u05d6\u05d9 \u05db\u05d9\u05d5\u05dd \u05d4\u05d5\u05d0 Sparse Dictionary Learning \u05d5\u05d1\u05e2\u05d9\u05e7\u05e8 Sparse Autoencoders. \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05d5\u05d0 \u05dc\u05d0\u05de\u05df \u05e8\u05e9\u05ea \u05e7\u05d8\u05e0\u05d4 \u05e9\"\u05ea\u05e4\u05e8\u05e9\" \u05d0\u05ea \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea \u05d4\u05d2\u05d3\u05d5\u05dc\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d1\u05e1\u05d9\u05e1 \u05d3\u05dc\u05d9\u05dc \u05e9\u05dc \"\u05ea\u05db\u05d5\u05e0\u05d5\u05ea\". \u05d0\u05dc\u05d5 \u05d4\u05dc\u05d9\u05d9\u05d8\u05e0\u05d8\u05d9\u05dd. \u05d0\u05da \u05d1\u05e4\u05d5\u05e2\u05dc, \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d0\u05de\u05e0\u05dd \u05de\u05d5\u05e6\u05d0\u05ea \u05db\u05d9\u05d5\u05d5\u05e0\u05d9\u05dd \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d9\u05dd, \u05d0\u05da \u05dc\u05d0 \u05de\u05e1\u05d1\u05d9\u05e8\u05d4 \u05d0\u05d9\u05da \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05e2\u05e6\u05de\u05d5 \u05de\u05ea\u05d1\u05e6\u05e2. \u05d4\u05dc\u05d9\u05d9\u05d8\u05e0\u05d8\u05d9\u05dd \u05d4\u05dd \u05ea\u05de\u05d5\u05e0\u05d4 \u05e1\u05d8\u05d8\u05d9\u05ea \u05e9\u05dc \"\u05de\u05d4 \u05d4\u05d5\u05e4\u05e2\u05dc\" \u05d5\u05dc\u05d0 \u05ea\u05d9\u05d0\u05d5\u05e8 \u05e9\u05dc \u05d4\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05e9\u05de\u05d9\u05d5\u05e9\u05dd.',\n  '\u05d9\u05e9 \u05d2\u05dd \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05d4\u05d5\u05ea\u05d9\u05d5\u05ea: \u05d4\u05e4\u05e2\u05e8 \u05d1\u05d9\u05df \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d4\u05d0\u05de\u05d9\u05ea\u05d9\u05d5\u05ea \u05dc\u05e9\u05d9\u05d7\u05d6\u05d5\u05e8\u05df \u05d2\u05d3\u05d5\u05dc. \u05d4\u05de\u05d9\u05d3\u05e2 \u05d4\u05d2\u05d0\u05d5\u05de\u05d8\u05e8\u05d9 \u05d1\u05d9\u05df \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d4\u05d5\u05dc\u05da \u05dc\u05d0\u05d9\u05d1\u05d5\u05d3. \u05d4\u05d4\u05e0\u05d7\u05d4 \u05e9\u05d4\u05db\u05d5\u05dc \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 \u05e8\u05d7\u05d5\u05e7\u05d4 \u05de\u05dc\u05d4\u05d9\u05d5\u05ea \u05e0\u05db\u05d5\u05e0\u05d4. \u05d5\u05d4\u05d2\u05e8\u05d5\u05e2 \u05de\u05db\u05dc \u05d4\u05d9\u05d0 \u05d4\u05e2\u05d5\u05d1\u05d3\u05d4 \u05e9\u05d0\u05d9\u05df \u05d1\u05db\u05dc\u05dc \u05ea\u05d9\u05d0\u05d5\u05e8\u05d9\u05d4 \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e9\u05de\u05e1\u05d1\u05d9\u05e8\u05d4 \u05de\u05d4\u05d9 \"\u05ea\u05db\u05d5\u05e0\u05d4\", \u05d0\u05d9\u05da \u05d4\u05d9\u05d0 \u05e0\u05d5\u05e6\u05e8\u05ea, \u05d5\u05de\u05d4 \u05d4\u05d5\u05e4\u05da \u05d0\u05d5\u05ea\u05d4 \u05dc\u05d9\u05d7\u05d9\u05d3\u05d4 \u05d1\u05e1\u05d9\u05e1\u05d9\u05ea \u05e9\u05dc \u05d4\u05d1\u05e0\u05d4.',\n  '\u05de\u05db\u05d0\u05df \u05e2\u05d5\u05dc\u05d4 \u05db\u05d9\u05d5\u05d5\u05df \u05d7\u05d3\u05e9\u05e0\u05d9: \u05d0\u05d5\u05dc\u05d9 \u05d4\u05d3\u05e8\u05da \u05d4\u05e0\u05db\u05d5\u05e0\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d0 \u05dc\u05e4\u05e8\u05e9 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05d7\u05e8\u05d9 \u05e9\u05d0\u05d5\u05de\u05e0\u05d5, \u05d0\u05dc\u05d0 \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05e8\u05d0\u05e9. \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e2\u05dd \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d3\u05d9\u05e1\u05e7\u05e8\u05d8\u05d9\u05d5\u05ea, \u05d0\u05db\u05d9\u05e4\u05ea \u05de\u05d5\u05d3\u05d5\u05dc\u05e8\u05d9\u05d5\u05ea, \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05e4\u05e2\u05dc\u05d4 \u05d3\u05dc\u05d9\u05dc\u05d5\u05ea \u05db\u05de\u05d5 Top-k \u05d0\u05d5 SoLU, \u05d0\u05d5 \u05de\u05d1\u05e0\u05d9\u05dd \u05db\u05de\u05d5 Mixture-of-Experts \u05e9\u05de\u05d7\u05dc\u05e7\u05d9\u05dd \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05dc\u05ea\u05ea-\u05de\u05d5\u05d3\u05d5\u05dc\u05d9\u05dd \u05d1\u05e8\u05d5\u05e8\u05d9\u05dd. \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d9\u05d9\u05e6\u05e8 \u05e8\u05e9\u05ea\u05d5\u05ea \u05e9\u05e0\u05d1\u05e0\u05d5\u05ea \"\u05d7\u05ea\u05d5\u05db\u05d5\u05ea \u05de\u05e8\u05d0\u05e9\" \u05e2\u05dd \u05e4\u05e8\u05e9\u05e0\u05d5\u05ea \u05dc\u05d0 \u05db\u05e0\u05d9\u05ea\u05d5\u05d7 \u05de\u05d0\u05d5\u05d7\u05e8 \u05d0\u05dc\u05d0 \u05db\u05d4\u05e0\u05d7\u05ea \u05d9\u05e1\u05d5\u05d3 \u05e9\u05dc \u05d4\u05d0\u05d9\u05de\u05d5\u05df.',\n  '\u05d2\u05dd \u05ea\u05d9\u05d0\u05d5\u05e8 \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e8\u05db\u05d9\u05d1 \u05d1\u05d5\u05d3\u05d3 \u05d4\u05d5\u05d0 \u05de\u05e9\u05d9\u05de\u05d4 \u05e7\u05e9\u05d4. \u05dc\u05de\u05e9\u05dc \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05de\u05e4\u05e2\u05d9\u05dc\u05d5\u05ea \u05d0\u05d5\u05ea\u05d5 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d1\u05dc\u05d1\u05dc\u05d5\u05ea. \u05e9\u05d9\u05d8\u05d5\u05ea \u05d9\u05d9\u05d7\u05d5\u05e1 \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05d1\u05e2\u05d9\u05d9\u05ea\u05d9\u05d5\u05ea \u05ea\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea \u05d5\u05e4\u05e8\u05e7\u05d8\u05d9\u05ea. \u05e1\u05d9\u05e0\u05ea\u05d6\u05ea \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d1\u05e0\u05d9\u05d9\u05ea \u05e7\u05dc\u05d8 \u05e9\u05de\u05e4\u05e2\u05d9\u05dc \u05e8\u05db\u05d9\u05d1 \u05e2\u05dc\u05d5\u05dc\u05d4 \u05dc\u05d9\u05d9\u05e6\u05e8 \u05d3\u05d9\u05de\u05d5\u05d9\u05d9\u05dd \u05dc\u05d0 \u05d0\u05d9\u05e0\u05e4\u05d5\u05e8\u05de\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05d4\u05de\u05d1\u05d8\u05d9\u05d7\u05d5\u05ea \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05df \u05d0\u05d9\u05e0\u05d8\u05e8\u05d5\u05d5\u05e0\u05e6\u05d9\u05d5\u05ea \u05e1\u05d9\u05d1\u05ea\u05d9\u05d5\u05ea: \u05e9\u05d9\u05e0\u05d5\u05d9 \u05e9\u05dc \u05e2\u05e8\u05da \u05e4\u05e0\u05d9\u05de\u05d9, \u05d5\u05d1\u05d7\u05d9\u05e0\u05d4 \u05e9\u05dc \u05d4\u05d4\u05e9\u05e4\u05e2\u05d4 \u05e2\u05dc \u05d4\u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05d7\u05d9\u05e6\u05d5\u05e0\u05d9\u05ea. \u05db\u05d0\u05df \u05e0\u05db\u05e0\u05e1\u05d9\u05dd \u05dc\u05ea\u05de\u05d5\u05e0\u05d4 \u05d2\u05dd steering \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d7\u05d3\u05e8\u05d4 \u05e9\u05dc \u05db\u05d9\u05d5\u05d5\u05df \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9 \u05dc\u05de\u05e8\u05d7\u05d1 \u05d4\u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d5\u05d2\u05dd \u05e9\u05d9\u05de\u05d5\u05e9 logit lens \u05db\u05d3\u05d9 \u05dc\u05e4\u05e2\u05e0\u05d7 \u05d4\u05e9\u05e4\u05e2\u05d4 \u05d9\u05e9\u05d9\u05e8\u05d4 \u05e2\u05dc \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05e2\u05dc \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d5\u05ea \u05d4\u05e8\u05e9\u05ea.',\n  '\u05d4\u05d1\u05e2\u05d9\u05d4 \u05d4\u05d2\u05d3\u05d5\u05dc\u05d4 \u05d4\u05d9\u05d0 \u05e9\u05d4\u05e8\u05d1\u05d4 \u05de\u05d4\u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05e0\u05e9\u05de\u05e2\u05d9\u05dd \u05de\u05e9\u05db\u05e0\u05e2\u05d9\u05dd \u05d0\u05da \u05dc\u05d0 \u05e2\u05d5\u05de\u05d3\u05d9\u05dd \u05d1\u05de\u05d1\u05d7\u05df. \u05d4\u05dd \u05dc\u05d0 \u05d7\u05d5\u05d6\u05d9\u05dd \u05e7\u05d5\u05e0\u05d8\u05e8\u05e4\u05e7\u05d8\u05d5\u05d0\u05dc\u05d9\u05dd(\u05dc\u05d0 \u05de\u05e6\u05dc\u05d9\u05d7\u05d9\u05dd \u05dc\u05e0\u05d1\u05d0 \u05de\u05d4 \u05d4\u05d9\u05d4 \u05e7\u05d5\u05e8\u05d4 \u05d0\u05d9\u05dc\u05d5 \u05de\u05e9\u05d4\u05d5 \u05d4\u05d9\u05d4 \u05e9\u05d5\u05e0\u05d4 \u05d1\u05ea\u05d5\u05da \u05d4\u05de\u05d5\u05d3\u05dc), \u05dc\u05d0 \u05e2\u05d5\u05d6\u05e8\u05d9\u05dd \u05dc\u05d0\u05d1\u05d7\u05df \u05db\u05e9\u05dc\u05d9 \u05de\u05d5\u05d3\u05dc, \u05dc\u05d0 \u05de\u05d0\u05e4\u05e9\u05e8\u05d9\u05dd \u05ea\u05d9\u05e7\u05d5\u05df \u05d0\u05d5 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d1\u05e4\u05d5\u05e2\u05dc. \u05dc\u05db\u05df \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e1\u05d8 \u05e9\u05dc\u05dd \u05e9\u05dc \u05d3\u05e8\u05db\u05d9 \u05d0\u05d9\u05de\u05d5\u05ea: \u05d4\u05d0\u05dd \u05d4\u05d4\u05e1\u05d1\u05e8 \u05d7\u05d5\u05d6\u05d4 \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d0\u05d7\u05e8\u05d9 ablation? \u05d4\u05d0\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc \u05e7\u05d8\u05df \u05e9\u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05dd \u05d4\u05d4\u05e1\u05d1\u05e8 \u05e0\u05db\u05d5\u05df? \u05d4\u05d0\u05dd \u05d4\u05dc\u05d9\u05d9\u05d8\u05e0\u05d8\u05d9\u05dd ( \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e4\u05e0\u05d9\u05de\u05d9\u05d9\u05dd \u05db\u05de\u05d5 \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d0\u05d5 \u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05d9\u05dd) \u05de\u05e1\u05d9\u05d9\u05e2\u05d9\u05dd \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d1\u05d8\u05d9\u05d7\u05d5\u05ea \u05db\u05de\u05d5 \u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d5\u05db\u05df \u05de\u05d6\u05d9\u05e7? \u05d4\u05d0\u05dd \u05e0\u05d5\u05db\u05dc \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05e9\u05e0\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05e0\u05d4\u05d2\u05d5\u05ea \u05d4\u05de\u05d5\u05d3\u05dc?',\n  '\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05d2\u05dd \u05dc\u05d9\u05e6\u05d5\u05e8 \"\u05d0\u05d5\u05e8\u05d2\u05e0\u05d9\u05d6\u05de\u05d9\u05dd \u05de\u05d5\u05d3\u05dc\u05d9\u05d9\u05dd\" \u05e9\u05d4\u05dd \u05de\u05d4\u05d5\u05d5\u05d9\u05dd \u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d8\u05e0\u05d5\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05d5\u05ea, \u05e2\u05dd \u05de\u05d1\u05e0\u05d4 \u05e4\u05ea\u05d5\u05d7, \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d0\u05de\u05df \u05e9\u05d5\u05d1 \u05d5\u05e9\u05d5\u05d1 \u05d5\u05dc\u05d1\u05d3\u05d5\u05e7 \u05e2\u05dc\u05d9\u05d4\u05df \u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05e8\u05e9\u05e0\u05d5\u05ea. \u05db\u05de\u05d5 \u05e9\u05d4\u05d1\u05d9\u05d5\u05dc\u05d5\u05d2\u05d9\u05d4 \u05d4\u05ea\u05e7\u05d3\u05de\u05d4 \u05d3\u05e8\u05da \u05e2\u05d1\u05d5\u05d3\u05d4 \u05e2\u05dc \u05ea\u05e1\u05d9\u05e1\u05e0\u05d9\u05ea, \u05db\u05da \u05ea\u05d7\u05d5\u05dd \u05d6\u05d4 \u05d6\u05e7\u05d5\u05e7 \u05dc\u05e8\u05e4\u05e8\u05e0\u05e1 \u05e7\u05d1\u05d5\u05e2. \u05d6\u05d4\u05d5 \u05db\u05dc\u05d9 \u05ea\u05e9\u05ea\u05d9\u05ea\u05d9 \u05d7\u05e1\u05e8 \u05db\u05d9\u05d5\u05dd.',\n  '\u05d4\u05d7\u05dc\u05e7 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d1\u05d4\u05d9\u05e8 \u05e9\u05de\u05db\u05e0\u05d9\u05d6\u05dd \u05d0\u05d9\u05e0\u05d5 \u05e2\u05e0\u05d9\u05d9\u05df \u05d8\u05db\u05e0\u05d9 \u05d1\u05dc\u05d1\u05d3. \u05d4\u05d5\u05d0 \u05e0\u05d5\u05d2\u05e2 \u05dc\u05de\u05d3\u05d9\u05e0\u05d9\u05d5\u05ea, \u05dc\u05e0\u05d9\u05d8\u05d5\u05e8, \u05dc\u05d1\u05d8\u05d9\u05d7\u05d5\u05ea, \u05d5\u05dc\u05e9\u05d0\u05dc\u05d5\u05ea \u05e4\u05d9\u05dc\u05d5\u05e1\u05d5\u05e4\u05d9\u05d5\u05ea: \u05de\u05d4 \u05e0\u05d7\u05e9\u05d1 \u05d4\u05e1\u05d1\u05e8 \u05d8\u05d5\u05d1? \u05d0\u05d9\u05da \u05d0\u05e4\u05e9\u05e8 \u05dc\u05d7\u05d1\u05e8 \u05d1\u05d9\u05df \u05d4\u05de\u05d1\u05e0\u05d9\u05dd \u05d4\u05de\u05d9\u05e7\u05e8\u05d5\u05e1\u05e7\u05d5\u05e4\u05d9\u05d9\u05dd \u05dc\u05ea\u05e4\u05e7\u05d5\u05d3 \u05d2\u05dc\u05d5\u05d1\u05dc\u05d9? \u05d0\u05d9\u05dc\u05d5 \u05e2\u05e7\u05e8\u05d5\u05e0\u05d5\u05ea \u05db\u05dc\u05dc\u05d9\u05d9\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05d7\u05dc\u05e5 \u05de\u05e8\u05e9\u05ea\u05d5\u05ea \u05e9\u05dc\u05de\u05d3\u05d5 \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05d8\u05d5\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05de\u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd?',\n  '\u05d1\u05e1\u05d9\u05db\u05d5\u05dd, \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05de\u05ea\u05d1\u05d9\u05d9\u05e9 \u05dc\u05d5\u05de\u05e8 \u05d0\u05ea \u05d4\u05d0\u05de\u05ea: \u05d0\u05d9\u05df \u05e2\u05d3\u05d9\u05d9\u05df \u05ea\u05d9\u05d0\u05d5\u05e8\u05d9\u05d4 \u05de\u05e1\u05e4\u05e7\u05ea \u05dc\u05e4\u05d9\u05e8\u05d5\u05e7 \u05e8\u05e9\u05ea\u05d5\u05ea. \u05d4\u05d4\u05e0\u05d7\u05d5\u05ea \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d5\u05ea \u05e9\u05d1\u05e8\u05d9\u05e8\u05d9\u05d5\u05ea. \u05d4\u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05dc\u05d0 \u05d7\u05d9\u05d5\u05ea \u05dc\u05d1\u05d3 \u05d0\u05dc\u05d0 \u05d1\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05e2\u05dc. \u05d4\u05e4\u05e8\u05e9\u05e0\u05d5\u05ea \u05d7\u05d9\u05d9\u05d1\u05ea \u05dc\u05e7\u05e9\u05d5\u05e8 \u05de\u05d1\u05e0\u05d4 \u05dc\u05ea\u05e4\u05e7\u05d5\u05d3. \u05d5\u05d4\u05d3\u05e8\u05da \u05e7\u05d3\u05d9\u05de\u05d4, \u05d0\u05d5\u05dc\u05d9, \u05e2\u05d5\u05d1\u05e8\u05ea \u05dc\u05d0 \u05d3\u05e8\u05da \u05e4\u05e2\u05e0\u05d5\u05d7 \u05d0\u05dc\u05d0 \u05d3\u05e8\u05da \u05f4\u05d3\u05d6\u05d9\u05d9\u05df\u05f4 \u05d7\u05d3\u05e9 \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea\u05d5\u05ea\u2026',\n)\nfor ($i = 0; $i -lt $middleTexts.Count; $i++) {\n  $d.Paragraphs.Item($i + 2).Range.Text = $middleTexts[$i]\n}\n\n# Paragraphs 14-20 (old numbering): delete entirely (7 paragraphs)\nfor ($i = 0; $i -lt 7; $i++) {\n  $d.Paragraphs.Item(14).Range.Delete()\n}\n\n# Last paragraph: replace URL text\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = 'https://arxiv.org/abs/2501.16496'\n\nWrite-Output \"done count=$($d.Paragraphs.Count)\""}
